$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 61, pushing existing rows 61-82 down to 62-83.
$ws.Rows.Item(61).Insert()

# The new row 61 shares every "static" (market/product) field with the rest
# of this block, so seed it from the template row directly above (row 60)
# and then overwrite the weekly observation columns with the new values.
for ($col = 1; $col -le 20; $col++) {
    $ws.Cells.Item(61, $col).Value2 = $ws.Cells.Item(60, $col).Value2
}

$ws.Range("D61").Value2 = 44606
$ws.Range("L61").Value2 = "Primera"
$ws.Range("M61").Value2 = 120
$ws.Range("N61").Value2 = 40000
$ws.Range("O61").Value2 = 42000
$ws.Range("P61").Value2 = 41000
$ws.Range("S61").Value2 = 2050

# Row 61's date needs the same date number-format as the other "Fecha"
# cells in this column.
$ws.Range("D61").NumberFormat = $ws.Range("D60").NumberFormat
